{"js": "// Update the answer key table: replace each three-digit division\n// problem's text with its newly generated counterpart.\n// Each old value is unique in the document, so a plain text search +\n// replace on the whole body is sufficient and safe.\n\nconst replacements = [\n  { old: \"100\u00f79=11, 1\", new: \"137\u00f75=27, 2\" },\n  { old: \"222\u00f78=27, 6\", new: \"818\u00f76=136, 2\" },\n  { old: \"554\u00f74=138, 2\", new: \"432\u00f72=216, 0\" },\n  { old: \"147\u00f74=36, 3\", new: \"188\u00f74=47, 0\" },\n  { old: \"467\u00f75=93, 2\", new: \"972\u00f79=108, 0\" },\n  { old: \"519\u00f72=259, 1\", new: \"398\u00f72=199, 0\" },\n  { old: \"102\u00f78=12, 6\", new: \"671\u00f72=335, 1\" },\n  { old: \"450\u00f72=225, 0\", new: \"550\u00f72=275, 0\" },\n  { old: \"431\u00f76=71, 5\", new: \"639\u00f75=127, 4\" },\n  { old: \"774\u00f76=129, 0\", new: \"830\u00f75=166, 0\" },\n  { old: \"271\u00f74=67, 3\", new: \"725\u00f73=241, 2\" },\n  { old: \"918\u00f74=229, 2\", new: \"712\u00f77=101, 5\" },\n  { old: \"514\u00f79=57, 1\", new: \"330\u00f78=41, 2\" },\n  { old: \"285\u00f73=95, 0\", new: \"106\u00f78=13, 2\" },\n  { old: \"880\u00f78=110, 0\", new: \"452\u00f78=56, 4\" },\n  { old: \"756\u00f76=126, 0\", new: \"693\u00f77=99, 0\" },\n  { old: \"843\u00f73=281, 0\", new: \"585\u00f79=65, 0\" },\n  { old: \"610\u00f74=152, 2\", new: \"861\u00f76=143, 3\" },\n  { old: \"853\u00f77=121, 6\", new: \"992\u00f74=248, 0\" },\n  { old: \"802\u00f76=133, 4\", new: \"817\u00f76=136, 1\" },\n  { old: \"484\u00f78=60, 4\", new: \"173\u00f76=28, 5\" },\n  { old: \"392\u00f74=98, 0\", new: \"958\u00f78=119, 6\" },\n  { old: \"182\u00f74=45, 2\", new: \"633\u00f74=158, 1\" },\n  { old: \"752\u00f74=188, 0\", new: \"568\u00f75=113, 3\" },\n  { old: \"775\u00f78=96, 7\", new: \"904\u00f76=150, 4\" },\n];\n\nconst body = context.document.body;\n\nfor (const { old, new: newText } of replacements) {\n  const results = body.search(old, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the answer key table: replace each three-digit division\n# problem's text with its newly generated counterpart.\n# Mapping is old-text -> new-text, applied as whole-field Find/Replace\n# across the whole document body (each old value is unique).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"100\u00f79=11, 1\";  new = \"137\u00f75=27, 2\"},\n    @{old = \"222\u00f78=27, 6\";  new = \"818\u00f76=136, 2\"},\n    @{old = \"554\u00f74=138, 2\"; new = \"432\u00f72=216, 0\"},\n    @{old = \"147\u00f74=36, 3\";  new = \"188\u00f74=47, 0\"},\n    @{old = \"467\u00f75=93, 2\";  new = \"972\u00f79=108, 0\"},\n    @{old = \"519\u00f72=259, 1\"; new = \"398\u00f72=199, 0\"},\n    @{old = \"102\u00f78=12, 6\";  new = \"671\u00f72=335, 1\"},\n    @{old = \"450\u00f72=225, 0\"; new = \"550\u00f72=275, 0\"},\n    @{old = \"431\u00f76=71, 5\";  new = \"639\u00f75=127, 4\"},\n    @{old = \"774\u00f76=129, 0\"; new = \"830\u00f75=166, 0\"},\n    @{old = \"271\u00f74=67, 3\";  new = \"725\u00f73=241, 2\"},\n    @{old = \"918\u00f74=229, 2\"; new = \"712\u00f77=101, 5\"},\n    @{old = \"514\u00f79=57, 1\";  new = \"330\u00f78=41, 2\"},\n    @{old = \"285\u00f73=95, 0\";  new = \"106\u00f78=13, 2\"},\n    @{old = \"880\u00f78=110, 0\"; new = \"452\u00f78=56, 4\"},\n    @{old = \"756\u00f76=126, 0\"; new = \"693\u00f77=99, 0\"},\n    @{old = \"843\u00f73=281, 0\"; new = \"585\u00f79=65, 0\"},\n    @{old = \"610\u00f74=152, 2\"; new = \"861\u00f76=143, 3\"},\n    @{old = \"853\u00f77=121, 6\"; new = \"992\u00f74=248, 0\"},\n    @{old = \"802\u00f76=133, 4\"; new = \"817\u00f76=136, 1\"},\n    @{old = \"484\u00f78=60, 4\";  new = \"173\u00f76=28, 5\"},\n    @{old = \"392\u00f74=98, 0\";  new = \"958\u00f78=119, 6\"},\n    @{old = \"182\u00f74=45, 2\";  new = \"633\u00f74=158, 1\"},\n    @{old = \"752\u00f74=188, 0\"; new = \"568\u00f75=113, 3\"},\n    @{old = \"775\u00f78=96, 7\";  new = \"904\u00f76=150, 4\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null\n}\n"}
